$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC_AddCustomerTEST_003: give rows 3-6 their own distinct usernames
# (previously every row duplicated row 2's "mngr464118"). This grows the
# shared-string table with four new unique values.
$ws.Range("A3").Value = "mngr464119"
$ws.Range("A4").Value = "mngr464120"
$ws.Range("A5").Value = "mngr464121"
$ws.Range("A6").Value = "mngr464122"
